$wb = $excel.ActiveWorkbook

# ===== Overview =====
$ws = $wb.Worksheets.Item('Overview')

# Update cell values
$ws.Range('A2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.md'
$ws.Range('B2').Value() = 'Handed back: in sync with en-US'
$ws.Range('C2').Value() = 'Handed back: in sync with en-US'
$ws.Range('D2').Value() = '2016-27-17 12:27:15'
$ws.Range('A3').Value() = '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md'
$ws.Range('B3').Value() = 'In Translation'
$ws.Range('C3').Value() = 'In Translation'
$ws.Range('D3').Value() = '2016-28-17 12:28:04'
$ws.Range('A4').Value() = '6d848c28-db13-4a46-a063-c5f816a6eee9.md'
$ws.Range('B4').Value() = 'In Translation'
$ws.Range('C4').Value() = 'In Translation'
$ws.Range('D4').Value() = '2016-28-17 12:28:04'
$ws.Range('A5').Value() = '41f88da1-c4d9-447d-8202-77217c99f257.md'
$ws.Range('B5').Value() = 'Ready for handoff'
$ws.Range('C5').Value() = 'Ready for handoff'
$ws.Range('D5').Value() = '2016-26-17 12:26:55'

# Rebuild hyperlinks (engine requires full delete+re-add; Delete() clears whole sheet)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/0025cdd25a59cb4667a78e3515d59606be9c7e3c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.md')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/41f88da1-c4d9-447d-8202-77217c99f257.md', '', '', '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md', '', '', '6d848c28-db13-4a46-a063-c5f816a6eee9.md')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/6d848c28-db13-4a46-a063-c5f816a6eee9.md', '', '', '41f88da1-c4d9-447d-8202-77217c99f257.md')

# ===== zh-cn =====
$ws = $wb.Worksheets.Item('zh-cn')

# Update cell values
$ws.Range('A2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.md'
$ws.Range('B2').Value() = '.md'
$ws.Range('C2').Value() = 'Handed back: in sync with en-US'
$ws.Range('D2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf'
$ws.Range('E2').Value() = '2016-03-17 12:27:12'
$ws.Range('F2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.md'
$ws.Range('G2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf'
$ws.Range('H2').Value() = '2016-03-17 12:27:28'
$ws.Range('I2').Value() = 'Include'
$ws.Range('A3').Value() = '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md'
$ws.Range('B3').Value() = '.md'
$ws.Range('C3').Value() = 'In Translation'
$ws.Range('D3').Value() = '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.zh-cn.xlf'
$ws.Range('E3').Value() = '2016-03-17 12:27:59'
$ws.Range('H3').Value() = '0001-01-01 00:00:00'
$ws.Range('I3').Value() = 'Include'
$ws.Range('A4').Value() = '6d848c28-db13-4a46-a063-c5f816a6eee9.md'
$ws.Range('B4').Value() = '.md'
$ws.Range('C4').Value() = 'In Translation'
$ws.Range('D4').Value() = '6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.zh-cn.xlf'
$ws.Range('E4').Value() = '2016-03-17 12:27:59'
$ws.Range('H4').Value() = '0001-01-01 00:00:00'
$ws.Range('I4').Value() = 'Include'
$ws.Range('A5').Value() = '41f88da1-c4d9-447d-8202-77217c99f257.md'
$ws.Range('B5').Value() = '.md'
$ws.Range('C5').Value() = 'Ready for handoff'
$ws.Range('D5').Value() = '41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.zh-cn.xlf'
$ws.Range('E5').Value() = '2016-03-17 12:26:52'
$ws.Range('H5').Value() = '0001-01-01 00:00:00'
$ws.Range('I5').Value() = 'Include'

# Rebuild hyperlinks (engine requires full delete+re-add; Delete() clears whole sheet)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/0025cdd25a59cb4667a78e3515d59606be9c7e3c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.md')
$ws.Hyperlinks.Add($ws.Range('B2'), 'https://github.com/OpenLocalizationTest/oltest/blob/0025cdd25a59cb4667a78e3515d59606be9c7e3c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/70027f86eb7a7ad69a21342b666347f730b39eac/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5dffa56c4dc31cd86b71b996d3c3f006bfcfb997/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/63b293bb596a0ad119bf53a61ec051810ab57c8d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/41f88da1-c4d9-447d-8202-77217c99f257.md', '', '', '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md')
$ws.Hyperlinks.Add($ws.Range('B3'), 'https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/41f88da1-c4d9-447d-8202-77217c99f257.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/22e58273805f33c2677647823a2b5a483ccfe9e9/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.zh-cn.xlf', '', '', '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md', '', '', '6d848c28-db13-4a46-a063-c5f816a6eee9.md')
$ws.Hyperlinks.Add($ws.Range('B4'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/093a276b0cd5ff68513929f3f6ded1be2c11e100/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.zh-cn.xlf', '', '', '6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.zh-cn.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/6d848c28-db13-4a46-a063-c5f816a6eee9.md', '', '', '41f88da1-c4d9-447d-8202-77217c99f257.md')
$ws.Hyperlinks.Add($ws.Range('B5'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/6d848c28-db13-4a46-a063-c5f816a6eee9.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/093a276b0cd5ff68513929f3f6ded1be2c11e100/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.zh-cn.xlf', '', '', '41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.zh-cn.xlf')

# ===== de-de =====
$ws = $wb.Worksheets.Item('de-de')

# Update cell values
$ws.Range('A2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.md'
$ws.Range('B2').Value() = '.md'
$ws.Range('C2').Value() = 'Handed back: in sync with en-US'
$ws.Range('D2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf'
$ws.Range('E2').Value() = '2016-03-17 12:27:15'
$ws.Range('F2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.md'
$ws.Range('G2').Value() = 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf'
$ws.Range('H2').Value() = '2016-03-17 12:27:34'
$ws.Range('I2').Value() = 'Include'
$ws.Range('A3').Value() = '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md'
$ws.Range('B3').Value() = '.md'
$ws.Range('C3').Value() = 'In Translation'
$ws.Range('D3').Value() = '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.de-de.xlf'
$ws.Range('E3').Value() = '2016-03-17 12:28:04'
$ws.Range('H3').Value() = '0001-01-01 00:00:00'
$ws.Range('I3').Value() = 'Include'
$ws.Range('A4').Value() = '6d848c28-db13-4a46-a063-c5f816a6eee9.md'
$ws.Range('B4').Value() = '.md'
$ws.Range('C4').Value() = 'In Translation'
$ws.Range('D4').Value() = '6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.de-de.xlf'
$ws.Range('E4').Value() = '2016-03-17 12:28:04'
$ws.Range('H4').Value() = '0001-01-01 00:00:00'
$ws.Range('I4').Value() = 'Include'
$ws.Range('A5').Value() = '41f88da1-c4d9-447d-8202-77217c99f257.md'
$ws.Range('B5').Value() = '.md'
$ws.Range('C5').Value() = 'Ready for handoff'
$ws.Range('D5').Value() = '41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.de-de.xlf'
$ws.Range('E5').Value() = '2016-03-17 12:26:55'
$ws.Range('H5').Value() = '0001-01-01 00:00:00'
$ws.Range('I5').Value() = 'Include'

# Rebuild hyperlinks (engine requires full delete+re-add; Delete() clears whole sheet)
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range('A2'), 'https://github.com/OpenLocalizationTest/oltest/blob/0025cdd25a59cb4667a78e3515d59606be9c7e3c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.md')
$ws.Hyperlinks.Add($ws.Range('B2'), 'https://github.com/OpenLocalizationTest/oltest/blob/0025cdd25a59cb4667a78e3515d59606be9c7e3c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D2'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/429ef1a6776e341f98ffc15018ccb6655c4b334f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('F2'), 'https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/38627967e5f3150198bb9b741f3b060fe40b838c/e2e/d2562679-66db-4d8d-b08b-7655ded4db5e.md', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.md')
$ws.Hyperlinks.Add($ws.Range('G2'), 'https://github.com/OpenLocalizationTestOrg/olhandback/blob/493d5c4783874726d2742ec2304bf0032c0d7862/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf', '', '', 'd2562679-66db-4d8d-b08b-7655ded4db5e.3a4de8036aab7508b025e0b5427f8e0059623ec4.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A3'), 'https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/41f88da1-c4d9-447d-8202-77217c99f257.md', '', '', '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md')
$ws.Hyperlinks.Add($ws.Range('B3'), 'https://github.com/OpenLocalizationTest/oltest/blob/f7b81d6d8fb8168467e114f8293c6f205ad3bb80/e2e/41f88da1-c4d9-447d-8202-77217c99f257.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D3'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4bfc791aa8658b10d0cdbbdadb83a60a04bac409/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.de-de.xlf', '', '', '5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A4'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md', '', '', '6d848c28-db13-4a46-a063-c5f816a6eee9.md')
$ws.Hyperlinks.Add($ws.Range('B4'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D4'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc129c7cb7657066205a4713fae51c0f64eb4c30/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5ef8e9bd-35be-40f3-9198-99f53ac9bf03.c1d0d23e8d521b14e9715f8e15b78c0287baee79.de-de.xlf', '', '', '6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.de-de.xlf')
$ws.Hyperlinks.Add($ws.Range('A5'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/6d848c28-db13-4a46-a063-c5f816a6eee9.md', '', '', '41f88da1-c4d9-447d-8202-77217c99f257.md')
$ws.Hyperlinks.Add($ws.Range('B5'), 'https://github.com/OpenLocalizationTest/oltest/blob/53c8037ae9403aa58c1775d16f1c9ba6f02e0a8c/e2e/6d848c28-db13-4a46-a063-c5f816a6eee9.md', '', '', '.md')
$ws.Hyperlinks.Add($ws.Range('D5'), 'https://github.com/OpenLocalizationTestOrg/olhandoff/blob/fc129c7cb7657066205a4713fae51c0f64eb4c30/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/6d848c28-db13-4a46-a063-c5f816a6eee9.2f10a57c10a13fd7336db648fce05066e75a61a0.de-de.xlf', '', '', '41f88da1-c4d9-447d-8202-77217c99f257.6e96362eaab6c05922d9ddb877256190515f9df0.de-de.xlf')
